$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column L: expand abbreviated position codes to full words
$posMap = @{ "r" = "right"; "b" = "center"; "y" = "left" }
for ($r = 2; $r -le 361; $r++) {
    $cell = $ws.Range("L$r")
    $code = $cell.Value2
    if ($posMap.ContainsKey($code)) {
        $cell.Value = $posMap[$code]
    }
}

# Replace stray "face//face_NN.jpg" stimulus filenames with "book//book_NN.jpg"
$ws.Range("D28").Value = "book//book_19.jpg"
$ws.Range("D34").Value = "book//book_01.jpg"
$ws.Range("D39").Value = "book//book_26.jpg"
$ws.Range("D92").Value = "book//book_24.jpg"
$ws.Range("D99").Value = "book//book_11.jpg"
$ws.Range("D103").Value = "book//book_31.jpg"
$ws.Range("D114").Value = "book//book_21.jpg"
$ws.Range("A122").Value = "book//book_13.jpg"
$ws.Range("A128").Value = "book//book_29.jpg"
$ws.Range("C131").Value = "book//book_29.jpg"
$ws.Range("A134").Value = "book//book_16.jpg"
$ws.Range("C138").Value = "book//book_16.jpg"
$ws.Range("A140").Value = "book//book_04.jpg"
$ws.Range("C142").Value = "book//book_04.jpg"
$ws.Range("C145").Value = "book//book_04.jpg"
$ws.Range("A146").Value = "book//book_20.jpg"
$ws.Range("C148").Value = "book//book_20.jpg"
$ws.Range("C150").Value = "book//book_20.jpg"
$ws.Range("A152").Value = "book//book_08.jpg"
$ws.Range("C157").Value = "book//book_08.jpg"
$ws.Range("A158").Value = "book//book_30.jpg"
$ws.Range("C161").Value = "book//book_30.jpg"
$ws.Range("A164").Value = "book//book_09.jpg"
$ws.Range("A170").Value = "book//book_19.jpg"
$ws.Range("C173").Value = "book//book_19.jpg"
$ws.Range("C175").Value = "book//book_19.jpg"
$ws.Range("A176").Value = "book//book_28.jpg"
$ws.Range("C178").Value = "book//book_28.jpg"
$ws.Range("C180").Value = "book//book_28.jpg"
$ws.Range("D184").Value = "book//book_11.jpg"
$ws.Range("D192").Value = "book//book_03.jpg"
$ws.Range("D219").Value = "book//book_38.jpg"
$ws.Range("D231").Value = "book//book_16.jpg"
$ws.Range("D241").Value = "book//book_35.jpg"
$ws.Range("A242").Value = "book//book_07.jpg"
$ws.Range("A248").Value = "book//book_21.jpg"
$ws.Range("A254").Value = "book//book_18.jpg"
$ws.Range("C259").Value = "book//book_18.jpg"
$ws.Range("A260").Value = "book//book_01.jpg"
$ws.Range("C262").Value = "book//book_01.jpg"
$ws.Range("C263").Value = "book//book_01.jpg"
$ws.Range("C264").Value = "book//book_01.jpg"
$ws.Range("A266").Value = "book//book_23.jpg"
$ws.Range("A272").Value = "book//book_38.jpg"
$ws.Range("A278").Value = "book//book_05.jpg"
$ws.Range("A284").Value = "book//book_27.jpg"
$ws.Range("C287").Value = "book//book_27.jpg"
$ws.Range("C288").Value = "book//book_27.jpg"
$ws.Range("A290").Value = "book//book_32.jpg"
$ws.Range("C292").Value = "book//book_32.jpg"
$ws.Range("C295").Value = "book//book_32.jpg"
$ws.Range("A296").Value = "book//book_26.jpg"
$ws.Range("D325").Value = "book//book_37.jpg"
$ws.Range("D330").Value = "book//book_32.jpg"
$ws.Range("D337").Value = "book//book_28.jpg"
$ws.Range("D353").Value = "book//book_12.jpg"
